# Update column G ("K") values on the active worksheet.
# These new values replace the previous Strike# based figures with
# regenerated K counts (std/mean recalculated, s_vals written upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = 1
    4  = 0
    5  = 2
    6  = 2
    7  = 2
    8  = 3
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 3
    18 = 1
    19 = 1
    21 = 2
    22 = 1
    23 = 1
    25 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
